# "Generate Report for Handoff"
#
# b6ad80c2-875e-4bc7-917d-52fdd310ebdb (already handed back & in sync)
# keeps its "Handed back: in sync with en-US" status and moves up to row 2,
# while 851bafc6-2f75-471d-abc7-2252955aef0d moves to row 3 and gets a new
# status "Ready for handoff" plus a fresh handoff timestamp, now that a new
# handoff report has been generated for it.

$wb = $excel.ActiveWorkbook

$uuidB = "b6ad80c2-875e-4bc7-917d-52fdd310ebdb"
$uuidA = "851bafc6-2f75-471d-abc7-2252955aef0d"

$mdB = "$uuidB.md"
$mdA = "$uuidA.md"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $mdB
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

$wsOverview.Range("A3").Value = $mdA
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $mdB
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $mdA
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlfB = "$uuidB.1a4b3cfe6ebd1aa2233dee01d0aa202a1bebc3b3.zh-cn.xlf"
$zhXlfA = "$uuidA.c170aa21613d147fd44592e71bfd192465f0c038.zh-cn.xlf"

# Row 2 now carries the b6ad80c2 entry (already handed back)
$wsZh.Range("A2").Value = $mdB
$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("C2").Value = $zhXlfB
$wsZh.Range("D2").Value = "2016-02-17 09:41:58"
$wsZh.Range("E2").Value = $mdB
$wsZh.Range("F2").Value = $zhXlfB
$wsZh.Range("G2").Value = "2016-02-17 09:42:54"
$wsZh.Range("H2").Value = "Include"

# Row 3 now carries the 851bafc6 entry, freshly handed off
$wsZh.Range("A3").Value = $mdA
$wsZh.Range("B3").Value = $statusReady
$wsZh.Range("C3").Value = $zhXlfA
$wsZh.Range("D3").Value = "2016-02-17 09:44:03"
$wsZh.Range("E3").Value = $mdA
$wsZh.Range("F3").Value = $zhXlfA
$wsZh.Range("G3").Value = "2016-02-17 09:42:54"
$wsZh.Range("H3").Value = "Include"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $mdB
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = $zhXlfB
    } elseif ($addr -eq '$E$2') {
        $h.TextToDisplay = $mdB
    } elseif ($addr -eq '$F$2') {
        $h.TextToDisplay = $zhXlfB
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $mdA
    } elseif ($addr -eq '$C$3') {
        $h.TextToDisplay = $zhXlfA
    } elseif ($addr -eq '$E$3') {
        $h.TextToDisplay = $mdA
    } elseif ($addr -eq '$F$3') {
        $h.TextToDisplay = $zhXlfA
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlfB = "$uuidB.1a4b3cfe6ebd1aa2233dee01d0aa202a1bebc3b3.de-de.xlf"
$deXlfA = "$uuidA.c170aa21613d147fd44592e71bfd192465f0c038.de-de.xlf"

# Row 2 now carries the b6ad80c2 entry (already handed back)
$wsDe.Range("A2").Value = $mdB
$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("C2").Value = $deXlfB
$wsDe.Range("D2").Value = "2016-02-17 09:42:10"
$wsDe.Range("E2").Value = $mdB
$wsDe.Range("F2").Value = $deXlfB
$wsDe.Range("G2").Value = "2016-02-17 09:43:14"
$wsDe.Range("H2").Value = "Include"

# Row 3 now carries the 851bafc6 entry, freshly handed off
$wsDe.Range("A3").Value = $mdA
$wsDe.Range("B3").Value = $statusReady
$wsDe.Range("C3").Value = $deXlfA
$wsDe.Range("D3").Value = "2016-02-17 09:44:19"
$wsDe.Range("E3").Value = $mdA
$wsDe.Range("F3").Value = $deXlfA
$wsDe.Range("G3").Value = "2016-02-17 09:43:14"
$wsDe.Range("H3").Value = "Include"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $mdB
    } elseif ($addr -eq '$C$2') {
        $h.TextToDisplay = $deXlfB
    } elseif ($addr -eq '$E$2') {
        $h.TextToDisplay = $mdB
    } elseif ($addr -eq '$F$2') {
        $h.TextToDisplay = $deXlfB
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $mdA
    } elseif ($addr -eq '$C$3') {
        $h.TextToDisplay = $deXlfA
    } elseif ($addr -eq '$E$3') {
        $h.TextToDisplay = $mdA
    } elseif ($addr -eq '$F$3') {
        $h.TextToDisplay = $deXlfA
    }
}
